$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.785.91"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "2.263.72"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.01"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.98"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.32"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.35"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.62"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "2.614.16"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.20"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "2.267.67"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "41.677.56"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.52"
$ws.Range("E20").Value = "  +6.32%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.90"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.69"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.00"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.39"
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.18"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.88"
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.37"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.94"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "2.024.39"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.28"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.37"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0278"
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("E47").Value = "  +15.43%  "
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.58"
$ws.Range("E51").Value = "  +4.52%  "
